# TICKET4769: template update to satisfy latest request
#
# Insert a new "Number of Rooms" / "%%PACKAGE_NUM_ROOMS%%" row into the
# package template, directly above the existing "Flex Minimum/Max" row,
# pushing it (and everything below it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Flex Minimum/Max" currently lives on row 20 - insert a new blank row
# there (Excel shifts row 20 and everything after it down to row 21+,
# copying formatting from the row above, same as a manual row-insert).
$ws.Rows.Item(20).Insert()

# Restore the standard row height / auto-height flag for the freshly
# inserted row (Insert() leaves it as a bare default row).
$ws.Rows.Item(20).RowHeight = 15

# Populate the new row with the new template fields.
$ws.Range("A20").Value = "Number of Rooms"
$ws.Range("B20").Value = "%%PACKAGE_NUM_ROOMS%%"
